# Backlog_7.xlsx edit:
#  - Column "C" (Semana) on both sheets changes from the text label
#    "Semana 07" to the plain number 7 for every populated data row.
#    This also causes the now-unused "Semana 07" shared string to be
#    dropped and every following shared string to shift down by one
#    slot (handled automatically on save).
#  - The active sheet changes from "ITI" back to "SPN", and each
#    sheet's saved selection is updated to reflect where the user
#    left off.

$wb = $excel.ActiveWorkbook

$wsSPN = $wb.Worksheets.Item("SPN")
$wsITI = $wb.Worksheets.Item("ITI")

# SPN: rows 2-13 hold data in column C ("Semana 07" -> 7)
for ($r = 2; $r -le 13; $r++) {
    $wsSPN.Range("C$r").Value = 7
}

# ITI: rows 2-29 hold data in column C ("Semana 07" -> 7)
for ($r = 2; $r -le 29; $r++) {
    $wsITI.Range("C$r").Value = 7
}

# Restore each sheet's saved selection.
$wsITI.Range("C2:C29").Select()
$wsSPN.Activate()
$wsSPN.Range("F18").Select()
